$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row content -------------------------------------------------
# Set values in an order that mirrors the target shared-strings table
# (0 Height, 1 Last Name of Child, 2 First Name of Child, 3 Address...,
#  4 Birthdate, 5 Weight, 6 Belongs to IP Group?, 7 Sex,
#  8 Taking Micronutrient Supplementation?, 9 Last Name of Parent/Guardian,
#  10 First Name of Parent/Guardian)
$ws.Range("J1").Value = "Height (cm)"
$ws.Range("D1").Value = "Last Name of Child"
$ws.Range("E1").Value = "First Name of Child"
$ws.Range("A1").Value = "Address or Location" + [char]10 + "of Child's Residence"
$ws.Range("G1").Value = "Birthdate" + [char]10 + "(YYYY-MM-DD)"
$ws.Range("K1").Value = "Weight (cm)"
$ws.Range("H1").Value = "Belongs to IP Group?" + [char]10 + "(Yes/No)"
$ws.Range("F1").Value = "Sex" + [char]10 + "(Male/Female)"
$ws.Range("I1").Value = "Taking Micronutrient Supplementation?" + [char]10 + "(Yes/No)"
$ws.Range("B1").Value = "Last Name" + [char]10 + "of Parent/Guardian"
$ws.Range("C1").Value = "First Name" + [char]10 + "of Parent/Guardian"

# --- Formatting -----------------------------------------------------------
# Give every header cell the same look the original headers already used
# (bold white font on green fill, centered) by copying the format from an
# existing header cell, then turn wrap text on for the cells that need it.
$ws.Range("D1").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1:C1").WrapText = $true
$ws.Range("F1:K1").WrapText = $true

$ws.Rows.Item(1).RowHeight = 45

# --- Column widths (best effort match to target pixel widths) -------------
$ws.Columns.Item(1).ColumnWidth = 43.251
$ws.Columns.Item(2).ColumnWidth = 34.917667
$ws.Columns.Item(3).ColumnWidth = 32.584333
$ws.Columns.Item(4).ColumnWidth = 24.751
$ws.Columns.Item(5).ColumnWidth = 25.917667
$ws.Columns.Item(7).ColumnWidth = 22.584333
$ws.Columns.Item(8).ColumnWidth = 22.584333
$ws.Columns.Item(9).ColumnWidth = 29.584333
$ws.Columns.Item(10).ColumnWidth = 12.251
$ws.Columns.Item(11).ColumnWidth = 12.584333

# --- Selection --------------------------------------------------------------
$ws.Range("H1").Select()
